$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled "Reminder" -> "Remainder" category text (row 26, column B / "type")
$ws.Range("B26").Value2 = "Remainder"

# Correct the "difficulty" values that were wrong, causing the question not to show in some sections
$ws.Range("D10").Value2 = 1
$ws.Range("D17").Value2 = 1

# Update the selection/view: select D10 (also resets the scrolled topLeftCell back to default)
$ws.Range("D10").Select()
